{"js": "// Replace each \"NN\u00d7NN=NNNN\" arithmetic answer in the table with its new value.\n// Old values are each unique in the document, so a body.search() per pair\n// (matchCase + exact string) is unambiguous.\nconst pairs = [\n  [\"16\u00d778=1248\", \"91\u00d718=1638\"],\n  [\"23\u00d793=2139\", \"43\u00d712=516\"],\n  [\"50\u00d738=1900\", \"57\u00d791=5187\"],\n  [\"14\u00d726=364\", \"48\u00d714=672\"],\n  [\"40\u00d762=2480\", \"79\u00d774=5846\"],\n  [\"68\u00d799=6732\", \"12\u00d764=768\"],\n  [\"98\u00d726=2548\", \"19\u00d794=1786\"],\n  [\"86\u00d762=5332\", \"93\u00d793=8649\"],\n  [\"56\u00d763=3528\", \"51\u00d713=663\"],\n  [\"40\u00d776=3040\", \"30\u00d742=1260\"],\n  [\"69\u00d742=2898\", \"58\u00d795=5510\"],\n  [\"11\u00d765=715\", \"58\u00d721=1218\"],\n  [\"60\u00d754=3240\", \"75\u00d761=4575\"],\n  [\"95\u00d769=6555\", \"69\u00d794=6486\"],\n  [\"15\u00d797=1455\", \"94\u00d778=7332\"],\n  [\"39\u00d717=663\", \"31\u00d723=713\"],\n  [\"42\u00d777=3234\", \"85\u00d745=3825\"],\n  [\"18\u00d798=1764\", \"36\u00d721=756\"],\n  [\"26\u00d795=2470\", \"34\u00d735=1190\"],\n  [\"90\u00d772=6480\", \"33\u00d795=3135\"],\n  [\"52\u00d712=624\", \"20\u00d745=900\"],\n  [\"63\u00d753=3339\", \"75\u00d733=2475\"],\n  [\"35\u00d762=2170\", \"82\u00d782=6724\"],\n  [\"92\u00d739=3588\", \"16\u00d718=288\"],\n  [\"64\u00d722=1408\", \"27\u00d727=729\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each \"NN\u00d7NN=NNNN\" arithmetic answer in the table with its new value.\n# Old values are each unique in the document, so Find/Replace with\n# MatchWholeWord off but an exact literal string is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"16\u00d778=1248\", \"91\u00d718=1638\"),\n    @(\"23\u00d793=2139\", \"43\u00d712=516\"),\n    @(\"50\u00d738=1900\", \"57\u00d791=5187\"),\n    @(\"14\u00d726=364\", \"48\u00d714=672\"),\n    @(\"40\u00d762=2480\", \"79\u00d774=5846\"),\n    @(\"68\u00d799=6732\", \"12\u00d764=768\"),\n    @(\"98\u00d726=2548\", \"19\u00d794=1786\"),\n    @(\"86\u00d762=5332\", \"93\u00d793=8649\"),\n    @(\"56\u00d763=3528\", \"51\u00d713=663\"),\n    @(\"40\u00d776=3040\", \"30\u00d742=1260\"),\n    @(\"69\u00d742=2898\", \"58\u00d795=5510\"),\n    @(\"11\u00d765=715\", \"58\u00d721=1218\"),\n    @(\"60\u00d754=3240\", \"75\u00d761=4575\"),\n    @(\"95\u00d769=6555\", \"69\u00d794=6486\"),\n    @(\"15\u00d797=1455\", \"94\u00d778=7332\"),\n    @(\"39\u00d717=663\", \"31\u00d723=713\"),\n    @(\"42\u00d777=3234\", \"85\u00d745=3825\"),\n    @(\"18\u00d798=1764\", \"36\u00d721=756\"),\n    @(\"26\u00d795=2470\", \"34\u00d735=1190\"),\n    @(\"90\u00d772=6480\", \"33\u00d795=3135\"),\n    @(\"52\u00d712=624\", \"20\u00d745=900\"),\n    @(\"63\u00d753=3339\", \"75\u00d733=2475\"),\n    @(\"35\u00d762=2170\", \"82\u00d782=6724\"),\n    @(\"92\u00d739=3588\", \"16\u00d718=288\"),\n    @(\"64\u00d722=1408\", \"27\u00d727=729\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $null, $null, $find.Forward, $find.Wrap, $null, $find.Replacement.Text, 2)  # wdReplaceAll\n}\n"}
